$wb = $excel.ActiveWorkbook

# --- Update scraped_at timestamps (column K) on the "snapshot" sheet ---
$ws = $wb.Worksheets.Item("snapshot")
$ws.Range("K2").Value = "2025-11-08T07:02:45.456129+00:00"
$ws.Range("K3").Value = "2025-11-08T07:02:45.456166+00:00"
$ws.Range("K4").Value = "2025-11-08T07:02:45.456190+00:00"
$ws.Range("K5").Value = "2025-11-08T07:02:48.318846+00:00"
$ws.Range("K6").Value = "2025-11-08T07:02:48.318877+00:00"
$ws.Range("K7").Value = "2025-11-08T07:02:48.318895+00:00"
$ws.Range("K8").Value = "2025-11-08T07:02:50.624461+00:00"
$ws.Range("K9").Value = "2025-11-08T07:02:52.914848+00:00"
$ws.Range("K10").Value = "2025-11-08T07:02:52.914877+00:00"
$ws.Range("K11").Value = "2025-11-08T07:02:52.914895+00:00"
$ws.Range("K12").Value = "2025-11-08T07:02:55.619058+00:00"
$ws.Range("K13").Value = "2025-11-08T07:02:55.619088+00:00"
$ws.Range("K14").Value = "2025-11-08T07:02:55.619105+00:00"
$ws.Range("K15").Value = "2025-11-08T07:02:55.619120+00:00"
$ws.Range("K16").Value = "2025-11-08T07:03:01.072475+00:00"
$ws.Range("K17").Value = "2025-11-08T07:03:03.407651+00:00"
$ws.Range("K18").Value = "2025-11-08T07:03:06.213073+00:00"
$ws.Range("K19").Value = "2025-11-08T07:03:06.213102+00:00"
$ws.Range("K20").Value = "2025-11-08T07:03:06.213120+00:00"
$ws.Range("K21").Value = "2025-11-08T07:03:08.908819+00:00"
$ws.Range("K22").Value = "2025-11-08T07:03:11.322693+00:00"
$ws.Range("K23").Value = "2025-11-08T07:03:11.322721+00:00"
$ws.Range("K24").Value = "2025-11-08T07:03:14.123250+00:00"
$ws.Range("K25").Value = "2025-11-08T07:03:14.123280+00:00"
$ws.Range("K26").Value = "2025-11-08T07:03:14.123298+00:00"
$ws.Range("K27").Value = "2025-11-08T07:03:16.419196+00:00"
$ws.Range("K28").Value = "2025-11-08T07:03:16.419227+00:00"
$ws.Range("K29").Value = "2025-11-08T07:03:16.419245+00:00"
$ws.Range("K30").Value = "2025-11-08T07:03:16.419260+00:00"
$ws.Range("K31").Value = "2025-11-08T07:03:16.419275+00:00"
$ws.Range("K32").Value = "2025-11-08T07:03:19.142429+00:00"
$ws.Range("K33").Value = "2025-11-08T07:03:21.918807+00:00"
$ws.Range("K34").Value = "2025-11-08T07:03:21.918837+00:00"
$ws.Range("K35").Value = "2025-11-08T07:03:21.918856+00:00"
$ws.Range("K36").Value = "2025-11-08T07:03:24.627984+00:00"
$ws.Range("K37").Value = "2025-11-08T07:03:24.628011+00:00"
$ws.Range("K38").Value = "2025-11-08T07:03:24.628029+00:00"
$ws.Range("K39").Value = "2025-11-08T07:03:26.939387+00:00"
$ws.Range("K40").Value = "2025-11-08T07:03:26.939417+00:00"
$ws.Range("K41").Value = "2025-11-08T07:03:26.939435+00:00"
$ws.Range("K42").Value = "2025-11-08T07:03:26.939451+00:00"
$ws.Range("K43").Value = "2025-11-08T07:03:26.939468+00:00"
$ws.Range("K44").Value = "2025-11-08T07:03:26.939484+00:00"
$ws.Range("K45").Value = "2025-11-08T07:03:29.320924+00:00"
$ws.Range("K46").Value = "2025-11-08T07:03:29.320954+00:00"
$ws.Range("K47").Value = "2025-11-08T07:03:34.459660+00:00"
$ws.Range("K48").Value = "2025-11-08T07:03:34.459688+00:00"
$ws.Range("K49").Value = "2025-11-08T07:03:36.846196+00:00"
$ws.Range("K50").Value = "2025-11-08T07:03:36.846224+00:00"

# --- Remove the stray row 2 from the "new_injured" sheet ---
$ws2 = $wb.Worksheets.Item("new_injured")
$ws2.Rows.Item(2).Delete()
